# Update cryptocurrency price/volume figures in the existing worksheet.
# Numeric-looking text values are prefixed with a leading apostrophe so
# Excel stores them as text (matching the original inlineStr cell type)
# rather than reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.020.12'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '3.521.21'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''592.15'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '''133.76'
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("D7").Value = '3.521.26'
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = '''0.124'
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("D12").Value = '''0.387'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '4.125.92'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").Value = '''27.66'
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").Value = '3.523.24'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '65.032.98'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").Value = '''10.18'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '''14.44'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").Value = '''5.71'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = '''392.27'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").Value = '''0.579'
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("D24").Value = '''74.95'
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = '3.663.95'
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '''0.0000112'
$ws.Range("E27").Value = '  -3.78%  '
$ws.Range("D28").Value = '''7.70'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +8.09%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").Value = '''8.36'
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '3.522.89'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '''24.12'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("E37").Value = '  +6.22%  '
$ws.Range("E38").Value = '  +3.13%  '
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("D40").Value = '''168.50'
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("D41").Value = '''0.0814'
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("D42").Value = '''0.822'
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("E43").Value = '  +5.25%  '
$ws.Range("D44").Value = '''25.87'
$ws.Range("E44").Value = '  -4.06%  '
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").Value = '2.414.90'
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("E51").Value = '  +5.66%  '
